$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date line
Replace-Text "2024-10-11 Friday" "2024-10-12 Saturday"

# Note: "409÷7=" is both a source value (index 5) and a target value of
# another rule (998÷3= -> 409÷7=). Replace the original "409÷7=" cell
# FIRST so the later write of a new "409÷7=" text is not re-matched.
Replace-Text "409÷7=" "269÷4="

Replace-Text "998÷3=" "409÷7="
Replace-Text "446÷8=" "642÷9="
Replace-Text "508÷8=" "465÷4="
Replace-Text "672÷6=" "680÷5="
Replace-Text "820÷6=" "955÷9="
Replace-Text "218÷5=" "970÷2="
Replace-Text "390÷4=" "335÷9="
Replace-Text "584÷2=" "127÷8="
Replace-Text "440÷9=" "279÷9="
Replace-Text "923÷5=" "529÷2="
Replace-Text "903÷5=" "103÷9="
Replace-Text "981÷3=" "769÷7="
Replace-Text "593÷9=" "703÷8="
Replace-Text "825÷4=" "681÷9="
Replace-Text "147÷9=" "251÷5="
Replace-Text "950÷6=" "861÷2="
Replace-Text "750÷8=" "274÷4="
Replace-Text "797÷2=" "905÷5="
Replace-Text "272÷9=" "684÷3="
Replace-Text "688÷9=" "965÷6="
Replace-Text "317÷6=" "348÷4="
Replace-Text "917÷6=" "174÷9="
Replace-Text "461÷7=" "527÷3="
Replace-Text "234÷5=" "196÷2="
